# Re-sync each market-report row to the corrected weekly snapshot.
# (The price sheet is re-derived weekly; rows keep their position but
#  pick up a different date/volume/price tuple from the new pull.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44390
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 2400
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2450
$ws.Range("P2").Value = 1225

# Row 3
$ws.Range("D3").Value = 44363
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = 2650
$ws.Range("P3").Value = 1325

# Row 4
$ws.Range("D4").Value = 44601
$ws.Range("J4").Value = 270
$ws.Range("K4").Value = 2200
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2350
$ws.Range("P4").Value = 1175

# Row 5
$ws.Range("D5").Value = 44427
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1400
$ws.Range("P5").Value = 700

# Row 6
$ws.Range("D6").Value = 45021
$ws.Range("J6").Value = 450
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2111
$ws.Range("P6").Value = 1056

# Row 7
$ws.Range("D7").Value = 44817
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 950
$ws.Range("P7").Value = 475

# Row 8
$ws.Range("D8").Value = 44438
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 950
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 975
$ws.Range("P8").Value = 488

# Row 10
$ws.Range("D10").Value = 44266
$ws.Range("K10").Value = 1700
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = 1750
$ws.Range("P10").Value = 875

# Row 11
$ws.Range("D11").Value = 44936
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3500
$ws.Range("M11").Value = 3357
$ws.Range("P11").Value = 1678

# Row 12
$ws.Range("D12").Value = 45001
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2250
$ws.Range("P12").Value = 1125

# Row 13
$ws.Range("D13").Value = 44253
$ws.Range("J13").Value = 250
$ws.Range("M13").Value = 1900
$ws.Range("P13").Value = 950

# Row 14
$ws.Range("D14").Value = 44544
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 950
$ws.Range("P14").Value = 475

# Row 15
$ws.Range("D15").Value = 45092
$ws.Range("K15").Value = 1800
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 1900
$ws.Range("P15").Value = 950

# Row 16
$ws.Range("D16").Value = 44243
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = 1250
$ws.Range("P16").Value = 625

# Row 17
$ws.Range("D17").Value = 44302
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 950
$ws.Range("P17").Value = 475

# Row 18
$ws.Range("D18").Value = 44468
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("P18").Value = 475

# Row 19
$ws.Range("D19").Value = 44435
$ws.Range("J19").Value = 300

# Row 20
$ws.Range("D20").Value = 44291
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 1800
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 1900
$ws.Range("P20").Value = 950

# Row 22
$ws.Range("D22").Value = 44229
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1900
$ws.Range("P22").Value = 950

# Row 23
$ws.Range("D23").Value = 44447
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 950
$ws.Range("P23").Value = 475

# Row 24
$ws.Range("D24").Value = 44525
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 1400
$ws.Range("M24").Value = 1450
$ws.Range("P24").Value = 725

# Row 25
$ws.Range("D25").Value = 44202
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 950

# Row 26
$ws.Range("D26").Value = 44540
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 950
$ws.Range("P26").Value = 475

# Row 27
$ws.Range("D27").Value = 44161
$ws.Range("J27").Value = 270

# Row 28
$ws.Range("D28").Value = 44392
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = 1900
$ws.Range("P28").Value = 950

# Row 29
$ws.Range("D29").Value = 45096
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = 950
$ws.Range("P29").Value = 475

# Row 31
$ws.Range("D31").Value = 44971
$ws.Range("J31").Value = 350
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 2800
$ws.Range("M31").Value = 2671
$ws.Range("P31").Value = 1336

# Row 32
$ws.Range("D32").Value = 44403
$ws.Range("J32").Value = 250

# Row 33
$ws.Range("D33").Value = 44572
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 1400
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = 1450
$ws.Range("P33").Value = 725

# Row 34
$ws.Range("D34").Value = 44726
$ws.Range("J34").Value = 250
$ws.Range("K34").Value = 2500
$ws.Range("L34").Value = 2800
$ws.Range("M34").Value = 2650
$ws.Range("P34").Value = 1325

# Row 35
$ws.Range("D35").Value = 45036
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 210
$ws.Range("K35").Value = 2300
$ws.Range("L35").Value = 2500
$ws.Range("M35").Value = 2443
$ws.Range("P35").Value = 1222

# Row 36
$ws.Range("D36").Value = 45062
$ws.Range("J36").Value = 380
$ws.Range("K36").Value = 1800
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 1895
$ws.Range("P36").Value = 948

# Row 37
$ws.Range("D37").Value = 44172
$ws.Range("J37").Value = 200
$ws.Range("K37").Value = 1300
$ws.Range("M37").Value = 1400
$ws.Range("P37").Value = 700

# Row 38
$ws.Range("D38").Value = 45055
$ws.Range("J38").Value = 350
$ws.Range("K38").Value = 2300
$ws.Range("L38").Value = 2500
$ws.Range("M38").Value = 2386
$ws.Range("P38").Value = 1193

# Row 39
$ws.Range("D39").Value = 45070
$ws.Range("J39").Value = 270
$ws.Range("K39").Value = 1000
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = 1250
$ws.Range("P39").Value = 625

# Row 40
$ws.Range("D40").Value = 44385
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 2400
$ws.Range("M40").Value = 2450
$ws.Range("P40").Value = 1225

# Row 41
$ws.Range("D41").Value = 44972
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 800
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = 943
$ws.Range("P41").Value = 472

# Row 43
$ws.Range("D43").Value = 44257
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 1400
$ws.Range("L43").Value = 1500
$ws.Range("M43").Value = 1450
$ws.Range("P43").Value = 725

# Row 44
$ws.Range("D44").Value = 44616
$ws.Range("J44").Value = 270
$ws.Range("K44").Value = 1300
$ws.Range("L44").Value = 1500
$ws.Range("M44").Value = 1400
$ws.Range("P44").Value = 700

# Row 45
$ws.Range("D45").Value = 45008
$ws.Range("J45").Value = 250
$ws.Range("M45").Value = 2200
$ws.Range("P45").Value = 1100

# Row 46
$ws.Range("D46").Value = 44365
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = 1900
$ws.Range("P46").Value = 950
